$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column D (Fecha, date serials) and column M (Volumen) values.
# The underlying data for rows 3-10 got permuted across rows; set the
# resulting target values directly (raw date serial numbers, as the
# cells already carry the date number format).

$ws.Range("D3").Value2 = 44306
$ws.Range("M3").Value2 = 80

$ws.Range("D4").Value2 = 44313
$ws.Range("M4").Value2 = 120

$ws.Range("D5").Value2 = 44316
$ws.Range("M5").Value2 = 120

$ws.Range("D6").Value2 = 44322
$ws.Range("M6").Value2 = 60

$ws.Range("D7").Value2 = 44323
$ws.Range("M7").Value2 = 80

$ws.Range("D8").Value2 = 44309
$ws.Range("M8").Value2 = 80
$ws.Range("Q8").Value = "$/caja 14 kilos granel"
$ws.Range("S8").Value2 = 821
$ws.Range("T8").Value2 = 14

$ws.Range("D9").Value2 = 44302
$ws.Range("M9").Value2 = 80

$ws.Range("D10").Value2 = 44327
$ws.Range("M10").Value2 = 60
$ws.Range("Q10").Value = "$/caja 10 kilos empedrada"
$ws.Range("S10").Value2 = 11500
$ws.Range("T10").Value2 = 1
